$wb = $excel.ActiveWorkbook

# Sheet: 展览 (exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 165
$ws1.Range("F5").Value = 6
$ws1.Range("F6").Value = 5873
$ws1.Range("G6").Value = 95
$ws1.Range("F9").Value = 3961
$ws1.Range("F16").Value = 115
$ws1.Range("F17").Value = 136
$ws1.Range("F18").Value = 661
$ws1.Range("F19").Value = 3961
$ws1.Range("F20").Value = 143
$ws1.Range("F22").Value = 5498
$ws1.Range("F23").Value = 444
$ws1.Range("F24").Value = 2168
$ws1.Range("F25").Value = 141
$ws1.Range("F27").Value = 8194
$ws1.Range("F32").Value = 184
$ws1.Range("F36").Value = 286
$ws1.Range("F37").Value = 262
$ws1.Range("F41").Value = 1189
$ws1.Range("F45").Value = 1370
$ws1.Range("F46").Value = 2183
$ws1.Range("F48").Value = 240
$ws1.Range("F49").Value = 1227

# Sheet: 本地生活 (local life)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 610
$ws3.Range("F3").Value = 789

# Sheet: 全部类型 (all types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 165
$ws4.Range("F5").Value = 610
$ws4.Range("F6").Value = 789
$ws4.Range("F7").Value = 5873
$ws4.Range("G7").Value = 95
$ws4.Range("F9").Value = 3961
$ws4.Range("F15").Value = 115
$ws4.Range("F18").Value = 661
$ws4.Range("F19").Value = 3961
$ws4.Range("F21").Value = 143
$ws4.Range("F23").Value = 5498
$ws4.Range("F24").Value = 444
$ws4.Range("F25").Value = 2168
$ws4.Range("F26").Value = 141
$ws4.Range("F28").Value = 8194
$ws4.Range("F32").Value = 184
$ws4.Range("F35").Value = 286
$ws4.Range("F36").Value = 262
$ws4.Range("F39").Value = 1189
$ws4.Range("F43").Value = 1370
$ws4.Range("F44").Value = 2183
$ws4.Range("F47").Value = 240
$ws4.Range("F49").Value = 1227
